$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous data row (row 33) onto the new row 35
$ws.Range("B33:C33").Copy()
$ws.Range("B35:C35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the new Q&A entry as row 35
$ws.Range("B35").Value = 43294
$ws.Range("C35").Value = "Można się odwoływać w metodzie Main do innych metod z innych klas instancja.metoda. Jak się można odwołać, nie będąc w metodzie Main do innej klasy?"

$ws.Range("B35:C35").RowHeight = 28.5

# Update the view to match the new content location
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C35").Select()
